$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master_Linking_Log")

# Force the new rows to text format temporarily so values such as dates
# ("2017-10-28") and plain numbers ("9", "7") are stored as text / shared
# strings rather than being auto-converted to numeric/date values, matching
# the original column formatting used by the rest of the sheet.
$newRows = $ws.Range("A6:H7")
$newRows.NumberFormat = "@"

$ws.Range("A6").Value = "j"
$ws.Range("B6").Value = "j"
$ws.Range("C6").Value = "2017-10-28"
$ws.Range("D6").Value = "22:01:58.011740"
$ws.Range("E6").Value = "j"
$ws.Range("F6").Value = "9"
$ws.Range("G6").Value = "m"
$ws.Range("H6").Value = "j"

$ws.Range("A7").Value = "a"
$ws.Range("B7").Value = "d"
$ws.Range("C7").Value = "2017-10-28"
$ws.Range("D7").Value = "12:30:00"
$ws.Range("E7").Value = "j"
$ws.Range("F7").Value = "7"
$ws.Range("G7").Value = "m"
$ws.Range("H7").Value = "j"

# Restore the default number format so the cells don't carry an explicit
# style reference (matching the unstyled cells used in rows 2-5).
$newRows.NumberFormat = "General"
